$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 105508726.4694378, -209018146.6045218),
    @(3, 132646577.4596441, -241504486.6946509),
    @(4, 116529554.7152026, -213442943.3238954),
    @(5, 106092782.3304012, -210929361.0420707),
    @(6, 105991829.6311245, -210946344.5005981),
    @(7, 122395754.4997257, -232183773.5343286),
    @(8, 106102412.26284, -210661424.2871895),
    @(9, 117040622.1587794, -215321766.9774082),
    @(10, 103861518.1174055, -208247002.9008795),
    @(11, 117134202.8675836, -214123845.5357292),
    @(12, 138056256.2293121, -238975883.2737637),
    @(13, 127395291.6096901, -231756331.7183017),
    @(14, 102692553.5673765, -204683855.8293287),
    @(15, 109438266.882265, -195287931.0485279),
    @(16, 138596224.7624223, -238383033.7740149),
    @(17, 130020051.2517149, -234012094.1658826),
    @(18, 134644140.8865899, -243949740.2878812),
    @(19, 139608841.0539081, -235476230.8749945),
    @(20, 105596981.9961531, -209263009.2431625),
    @(21, 141599423.6154124, -238858232.1724758),
    @(22, 134875384.7845174, -229901430.6513903),
    @(23, 105594822.7153038, -210140147.9385734),
    @(24, 105845394.8503075, -210537490.5371626),
    @(25, 134781258.9465282, -235788498.9133089),
    @(26, 125182348.8449171, -226212495.6601892),
    @(27, 105460556.4692007, -209469655.8274356),
    @(28, 127417554.1764413, -218170442.130488),
    @(29, 110635064.8260035, -200693955.3049799),
    @(30, 136104466.1075085, -232184709.8015506),
    @(31, 134783690.3161738, -235791494.57205),
    @(32, 105451235.8885704, -209736329.8344363),
    @(33, 133553287.5962651, -243785214.8407176),
    @(34, 105716464.7564075, -210474208.3568346),
    @(35, 139066872.48509, -234520873.5788515),
    @(36, 138187606.9388872, -238041235.5106107),
    @(37, 111701037.9983043, -202577032.8259242),
    @(38, 113236264.3979196, -203448141.3676558),
    @(39, 102773319.0900762, -207250826.245541),
    @(40, 105017672.2254655, -208242418.4061436),
    @(41, 134701314.1507206, -235598052.3902366),
    @(42, 104439268.9421774, -208685705.4748633),
    @(43, 104063185.0347636, -206174733.721487),
    @(44, 116381563.0254742, -213354002.2535594),
    @(45, 100167826.3336622, -200491735.4238619),
    @(46, 106181608.3227286, -211176252.3923558),
    @(47, 104367728.3348256, -209373254.1300425),
    @(48, 104424001.3408905, -206884594.4189475),
    @(49, 137283281.5320168, -236821997.8006023),
    @(50, 106000975.046599, -209794092.8900959),
    @(51, 117159326.3814496, -203346207.1822138),
    @(52, 123465584.6433238, -233430269.622158),
    @(53, 130665300.368854, -238428189.4724956),
    @(54, 100367178.8657337, -200478489.6620883),
    @(55, 105352880.5020434, -208877513.5332501),
    @(56, 129647122.8856552, -233086310.6818279),
    @(57, 99652515.60497172, -202793563.913616),
    @(58, 124626453.9291485, -234583616.9577389),
    @(59, 112672866.8429443, -202306174.8753476),
    @(60, 120818759.6407465, -228976588.5435393),
    @(61, 136059248.6204675, -251416086.562739),
    @(62, 137281823.6174914, -236820256.6947784),
    @(63, 111557958.197442, -198188443.8184662),
    @(64, 122342850.5698752, -231739329.838951),
    @(65, 113504335.765346, -221679853.6071709),
    @(66, 117746875.3597333, -203145161.1829448),
    @(67, 123466691.2825127, -233431851.2426446),
    @(68, 110571376.8817002, -199724522.8641171),
    @(69, 110055601.0016471, -192394145.9747593),
    @(70, 115877033.3297141, -206198179.3881435),
    @(71, 130111079.5130356, -227255961.6400539),
    @(72, 105598637.6146314, -209265622.4662407),
    @(73, 134781913.8277764, -235789305.7885106),
    @(74, 120457917.0446891, -222901472.9667384),
    @(75, 104656956.3492973, -208998483.7739748),
    @(76, 140447377.957341, -241190465.6660586),
    @(77, 127382852.6466472, -223616842.4244619),
    @(78, 106918502.4314488, -193926183.2665913),
    @(79, 99678515.17212091, -202869844.3416095),
    @(80, 124630512.6759802, -234589366.480262),
    @(81, 122906404.1367475, -232588847.2752898),
    @(82, 123834612.5745275, -234081588.0408616),
    @(83, 125024745.8866452, -215332164.6101156),
    @(84, 114770796.1817679, -199968539.1906832),
    @(85, 136399735.7707716, -235342743.6443977),
    @(86, 124134250.592144, -234550627.4243446),
    @(87, 109970644.6036005, -195827838.7745652),
    @(88, 118735399.5336341, -225639454.349205),
    @(89, 99948772.93257391, -203348593.3003232),
    @(90, 113568385.9032365, -221247347.6185689),
    @(91, 99680708.61797616, -202873459.8895674),
    @(92, 110142394.5216578, -208562628.9073775),
    @(93, 99031260.30550426, -201602192.2960845),
    @(94, 99679793.69599606, -202871951.7895998),
    @(95, 99085419.68022577, -201758063.0948648),
    @(96, 117016909.536118, -209587252.9708545),
    @(97, 102545355.9689802, -206606667.7502118),
    @(98, 117672836.3618993, -216358649.6340961),
    @(99, 123978034.8062271, -234452325.5146987),
    @(100, 120470081.3920262, -228674520.3612784),
    @(101, 137093107.2131552, -232680201.6958987),
    @(102, 105511964.562603, -209023256.3673792),
    @(103, 114333012.3162275, -221797612.6148877),
    @(104, 117585726.1880774, -222315825.1105521),
    @(105, 107604625.910323, -209515713.5931183),
    @(106, 145569466.401966, -261425231.6330547),
    @(107, 130113878.5378258, -230877998.5930108),
    @(108, 105080389.1113953, -209884292.3578205),
    @(109, 120489265.4427217, -208929469.54725),
    @(110, 98461024.83522852, -197396514.180316),
    @(111, 115607022.0688982, -205529791.2928719),
    @(112, 124822083.9285265, -236031568.2958097),
    @(113, 100156661.0837503, -200172776.0156584),
    @(114, 100001447.5327858, -200022630.0171967),
    @(115, 141041486.155849, -256629137.1601039),
    @(116, 115944493.4328775, -208300507.8710816),
    @(117, 116477400.184276, -213021100.7277543),
    @(118, 134073855.6473529, -232856329.5791467),
    @(119, 117860571.0843874, -220944422.1572055),
    @(120, 103446922.0606982, -205335959.0468844),
    @(121, 99679180.36433366, -202870940.8090142),
    @(122, 99021577.51894188, -198404658.9435754),
    @(123, 99706732.31271403, -199480058.8985869),
    @(124, 124840981.4588256, -221515326.7091922),
    @(125, 101482637.130122, -201593533.2296944),
    @(126, 115547775.5122885, -211230821.352387),
    @(127, 109571273.3401862, -191719277.4940407),
    @(128, 99678284.8839162, -202869464.7461046),
    @(129, 116896498.8276936, -215237703.2306432),
    @(130, 101486613.3077085, -195192346.6161174),
    @(131, 111206429.3717718, -209942642.4337206),
    @(132, 100690301.4471169, -200528169.7813111),
    @(133, 99441616.50272426, -202479162.9386826),
    @(134, 113308239.1278594, -216662236.3820526),
    @(135, 117210532.9378552, -202388640.3289137),
    @(136, 104803634.4021675, -207065422.8186131),
    @(137, 115088902.8054884, -211475437.8546615),
    @(138, 108953363.5511799, -208830389.1402352),
    @(139, 132227328.7373793, -246882657.1680205),
    @(140, 120645077.3307229, -226992646.1945517),
    @(141, 100370172.5250075, -200483310.2279336),
    @(142, 121803530.9674991, -214929262.2081245),
    @(143, 115093850.2074087, -210639130.9912758),
    @(144, 103822953.4520017, -189930731.702317),
    @(145, 132877443.62134, -224023748.5784365),
    @(146, 101015466.2542255, -201723206.6404476),
    @(147, 119228449.4292004, -226634202.4660612),
    @(148, 140082692.4468106, -234409595.1433657),
    @(149, 119395384.9325225, -224644429.6894445),
    @(150, 107077921.004971, -210415575.5129712),
    @(151, 116599761.8003678, -215056264.1719109),
    @(152, 143874749.0634952, -258943761.0773299),
    @(153, 122093307.7635642, -229543200.8580566),
    @(154, 105590085.9079321, -208653782.3603525),
    @(155, 132597887.1665408, -232357831.3243726),
    @(156, 103470805.8748131, -204822101.9525426),
    @(157, 127180613.9277396, -221078511.4456128),
    @(158, 116386930.9560903, -202553336.0838474),
    @(159, 101583352.7685994, -202736787.3699743),
    @(160, 139999211.642137, -234218738.9120067),
    @(161, 129403918.6188121, -220561130.1128934),
    @(162, 118958854.6132308, -225924513.2256453),
    @(163, 127361799.8553083, -222684836.8283646),
    @(164, 105214205.3888756, -208202346.972273),
    @(165, 130001514.6469524, -221096251.1072609),
    @(166, 137717459.0523186, -235681984.5561872),
    @(167, 105128541.7492181, -207962133.4922707),
    @(168, 121060300.4620544, -220297557.7979507),
    @(169, 136483512.4430443, -226734556.8126839),
    @(170, 124548501.1505335, -217902651.7519901),
    @(171, 128464498.0296568, -219784764.8501419),
    @(172, 136901542.6570788, -236224593.5605564),
    @(173, 140229425.9082486, -234134835.1599286),
    @(174, 115411245.186507, -213155727.6359349),
    @(175, 137088898.2118142, -228123017.841619),
    @(176, 138936348.0628513, -239502436.3370571),
    @(177, 116306647.5220219, -212590809.3727053),
    @(178, 120195467.8391005, -227955826.4722086),
    @(179, 119185180.0915828, -226213221.2188751),
    @(180, 132649677.9409595, -241194673.5050797)
)

foreach ($item in $data) {
    $r = $item[0]
    $b = $item[1]
    $c = $item[2]
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = 48
    $ws.Cells.Item($r, 5).Value = "U"
}

Write-Host "Done applying changes to" $data.Count "rows"